# Clean and organize the interview-schedule workbook:
#  - Delete the two leftover candidates (Ashlyn/Reich, Aurelio/Molina) from
#    the "Facebook" sheet, and the two blank trailing rows on "Apple".
#  - Normalize both sheets' header rows to fname / lname / Time.
#  - Make "Facebook" the active tab, with A2:C3 selected.

$wb = $excel.ActiveWorkbook

$apple = $wb.Worksheets.Item("Apple")
$facebook = $wb.Worksheets.Item("Facebook")

# --- Facebook: drop the two rows that aren't part of the schedule anymore ---
$facebook.Rows.Item(2).Delete()
$facebook.Rows.Item(2).Delete()

# --- Apple: drop the two trailing blank rows ---
$apple.Rows.Item(7).Delete()
$apple.Rows.Item(7).Delete()

# --- Normalize headers on both sheets ---
$apple.Range("A1").Value = "fname"
$apple.Range("B1").Value = "lname"
$apple.Range("C1").Value = "Time"

$facebook.Range("A1").Value = "fname"
$facebook.Range("B1").Value = "lname"
$facebook.Range("C1").Value = "Time"

# --- Selection / active tab ---
$apple.Range("A2:C3").Select()

$facebook.Activate()
$facebook.Range("A2:C3").Select()
